$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UserTestData")

# New test data row
$ws.Range("A4").Value = "Test Edit User details API"
$ws.Range("B4").Value = 8106977

# Update existing UserId value in row 2
$ws.Range("B2").Value = 8107043

# D4 should carry the Hyperlink cell style (same as D3) but remain empty
$ws.Range("D4").Style = $ws.Range("D3").Style

# Widen columns A, C, D (re-"best fit") to account for the new, wider content
$ws.Columns.Item(1).ColumnWidth = 20.65
$ws.Columns.Item(3).ColumnWidth = 18.65
$ws.Columns.Item(4).ColumnWidth = 20.65

# Match the saved selection/active cell from the authored workbook
$ws.Range("D12").Select()
